# Automatic update of files.
# Update the "Förändrad" (Changed) date column (C) for rows 2-54
# from 2023-09-15 (serial 45184) to 2023-09-16 (serial 45185).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 54; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
